# Update the "Förändrad" (Changed) date column (C) for all data rows
# (rows 2-524) from 2023-09-10 (serial 45179) to 2023-09-11 (serial 45180).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C524").Value = 45180
